$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the columns that are going away (First Name, Last Name, Time In, Time Out) ---
# Keep old C (Email), F (Date Time In...), G (Date Time Out...) which become the new A, B, C.
$ws.Columns("A").Delete()       # First Name gone; Last Name -> A, Email -> B, ...
$ws.Columns("A").Delete()       # Last Name gone; Email -> A, Time In -> B, ...
$ws.Columns("B:C").Delete()     # Time In / Time Out gone; Date Time In -> B, Date Time Out -> C

# --- Re-purpose the headers for the new layout ---
$ws.Range("A1").Value = "Emplyee Number"
$ws.Range("B1").Value = "Date Time In(YYYY-MM-DD HH:MM:SS)"
$ws.Range("C1").Value = "Date Time Out(YYY-MM-DD HH:MM:SS)"

# --- Trim the sheet from 14 data rows down to 12 ---
$ws.Rows("13:14").Delete()

# --- Widen the new last column ---
$ws.Columns("C").ColumnWidth = 36.1

# --- Store everything as text (matches the new free-form "HH:MM:SS" headers) ---
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:A12").NumberFormat = "@"

# --- Match the saved selection ---
$ws.Range("B8").Select() | Out-Null

Write-Host "done"
